# fix(import): add antenne column
# Adds a new "antenne" column (K) with header + sample value "MONTREUIL",
# and moves the active selection to the next empty row/column (K3) as a
# hint for where the next value would be entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "antenne"
$ws.Range("K2").Value = "MONTREUIL"

$ws.Range("K3").Select()
